$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite rows 27-91 (and extend dimension) per the updated dataset.
# Rows 2-26 (Brasil 2016-2022) are unchanged.

$ws.Cells.Item(27, 1).Value = "Brasil"
$ws.Cells.Item(27, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(27, 3).Value = "31/12/2023"
$ws.Cells.Item(27, 4).Value = 62.3

$ws.Cells.Item(28, 1).Value = "Brasil"
$ws.Cells.Item(28, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(28, 3).Value = "31/12/2023"
$ws.Cells.Item(28, 4).Value = 6

$ws.Cells.Item(29, 1).Value = "Brasil"
$ws.Cells.Item(29, 2).Value = "Alugado"
$ws.Cells.Item(29, 3).Value = "31/12/2023"
$ws.Cells.Item(29, 4).Value = 22.4

$ws.Cells.Item(30, 1).Value = "Brasil"
$ws.Cells.Item(30, 2).Value = "Cedido"
$ws.Cells.Item(30, 3).Value = "31/12/2023"
$ws.Cells.Item(30, 4).Value = 9

$ws.Cells.Item(31, 1).Value = "Brasil"
$ws.Cells.Item(31, 2).Value = "Outra condição"
$ws.Cells.Item(31, 3).Value = "31/12/2023"
$ws.Cells.Item(31, 4).Value = 0.3

$ws.Cells.Item(32, 1).Value = "Nordeste"
$ws.Cells.Item(32, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(32, 3).Value = "31/12/2016"
$ws.Cells.Item(32, 4).Value = 73

$ws.Cells.Item(33, 1).Value = "Nordeste"
$ws.Cells.Item(33, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(33, 3).Value = "31/12/2016"
$ws.Cells.Item(33, 4).Value = 3.6

$ws.Cells.Item(34, 1).Value = "Nordeste"
$ws.Cells.Item(34, 2).Value = "Alugado"
$ws.Cells.Item(34, 3).Value = "31/12/2016"
$ws.Cells.Item(34, 4).Value = 15

$ws.Cells.Item(35, 1).Value = "Nordeste"
$ws.Cells.Item(35, 2).Value = "Cedido"
$ws.Cells.Item(35, 3).Value = "31/12/2016"
$ws.Cells.Item(35, 4).Value = 8.199999999999999

$ws.Cells.Item(36, 1).Value = "Nordeste"
$ws.Cells.Item(36, 2).Value = "Outra condição"
$ws.Cells.Item(36, 3).Value = "31/12/2016"
$ws.Cells.Item(36, 4).Value = 0.2

$ws.Cells.Item(37, 1).Value = "Nordeste"
$ws.Cells.Item(37, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(37, 3).Value = "31/12/2017"
$ws.Cells.Item(37, 4).Value = 73.09999999999999

$ws.Cells.Item(38, 1).Value = "Nordeste"
$ws.Cells.Item(38, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(38, 3).Value = "31/12/2017"
$ws.Cells.Item(38, 4).Value = 3.5

$ws.Cells.Item(39, 1).Value = "Nordeste"
$ws.Cells.Item(39, 2).Value = "Alugado"
$ws.Cells.Item(39, 3).Value = "31/12/2017"
$ws.Cells.Item(39, 4).Value = 15

$ws.Cells.Item(40, 1).Value = "Nordeste"
$ws.Cells.Item(40, 2).Value = "Cedido"
$ws.Cells.Item(40, 3).Value = "31/12/2017"
$ws.Cells.Item(40, 4).Value = 8.199999999999999

$ws.Cells.Item(41, 1).Value = "Nordeste"
$ws.Cells.Item(41, 2).Value = "Outra condição"
$ws.Cells.Item(41, 3).Value = "31/12/2017"
$ws.Cells.Item(41, 4).Value = 0.2

$ws.Cells.Item(42, 1).Value = "Nordeste"
$ws.Cells.Item(42, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(42, 3).Value = "31/12/2018"
$ws.Cells.Item(42, 4).Value = 71.5

$ws.Cells.Item(43, 1).Value = "Nordeste"
$ws.Cells.Item(43, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(43, 3).Value = "31/12/2018"
$ws.Cells.Item(43, 4).Value = 3.7

$ws.Cells.Item(44, 1).Value = "Nordeste"
$ws.Cells.Item(44, 2).Value = "Alugado"
$ws.Cells.Item(44, 3).Value = "31/12/2018"
$ws.Cells.Item(44, 4).Value = 15.7

$ws.Cells.Item(45, 1).Value = "Nordeste"
$ws.Cells.Item(45, 2).Value = "Cedido"
$ws.Cells.Item(45, 3).Value = "31/12/2018"
$ws.Cells.Item(45, 4).Value = 9

$ws.Cells.Item(46, 1).Value = "Nordeste"
$ws.Cells.Item(46, 2).Value = "Outra condição"
$ws.Cells.Item(46, 3).Value = "31/12/2018"
$ws.Cells.Item(46, 4).Value = 0.2

$ws.Cells.Item(47, 1).Value = "Nordeste"
$ws.Cells.Item(47, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(47, 3).Value = "31/12/2019"
$ws.Cells.Item(47, 4).Value = 72.40000000000001

$ws.Cells.Item(48, 1).Value = "Nordeste"
$ws.Cells.Item(48, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(48, 3).Value = "31/12/2019"
$ws.Cells.Item(48, 4).Value = 3.2

$ws.Cells.Item(49, 1).Value = "Nordeste"
$ws.Cells.Item(49, 2).Value = "Alugado"
$ws.Cells.Item(49, 3).Value = "31/12/2019"
$ws.Cells.Item(49, 4).Value = 15.7

$ws.Cells.Item(50, 1).Value = "Nordeste"
$ws.Cells.Item(50, 2).Value = "Cedido"
$ws.Cells.Item(50, 3).Value = "31/12/2019"
$ws.Cells.Item(50, 4).Value = 8.6

$ws.Cells.Item(51, 1).Value = "Nordeste"
$ws.Cells.Item(51, 2).Value = "Outra condição"
$ws.Cells.Item(51, 3).Value = "31/12/2019"
$ws.Cells.Item(51, 4).Value = 0.2

$ws.Cells.Item(52, 1).Value = "Nordeste"
$ws.Cells.Item(52, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(52, 3).Value = "31/12/2022"
$ws.Cells.Item(52, 4).Value = 71

$ws.Cells.Item(53, 1).Value = "Nordeste"
$ws.Cells.Item(53, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(53, 3).Value = "31/12/2022"
$ws.Cells.Item(53, 4).Value = 3.2

$ws.Cells.Item(54, 1).Value = "Nordeste"
$ws.Cells.Item(54, 2).Value = "Alugado"
$ws.Cells.Item(54, 3).Value = "31/12/2022"
$ws.Cells.Item(54, 4).Value = 17.2

$ws.Cells.Item(55, 1).Value = "Nordeste"
$ws.Cells.Item(55, 2).Value = "Cedido"
$ws.Cells.Item(55, 3).Value = "31/12/2022"
$ws.Cells.Item(55, 4).Value = 8.4

$ws.Cells.Item(56, 1).Value = "Nordeste"
$ws.Cells.Item(56, 2).Value = "Outra condição"
$ws.Cells.Item(56, 3).Value = "31/12/2022"
$ws.Cells.Item(56, 4).Value = 0.2

$ws.Cells.Item(57, 1).Value = "Nordeste"
$ws.Cells.Item(57, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(57, 3).Value = "31/12/2023"
$ws.Cells.Item(57, 4).Value = 70.40000000000001

$ws.Cells.Item(58, 1).Value = "Nordeste"
$ws.Cells.Item(58, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(58, 3).Value = "31/12/2023"
$ws.Cells.Item(58, 4).Value = 2.9

$ws.Cells.Item(59, 1).Value = "Nordeste"
$ws.Cells.Item(59, 2).Value = "Alugado"
$ws.Cells.Item(59, 3).Value = "31/12/2023"
$ws.Cells.Item(59, 4).Value = 17.9

$ws.Cells.Item(60, 1).Value = "Nordeste"
$ws.Cells.Item(60, 2).Value = "Cedido"
$ws.Cells.Item(60, 3).Value = "31/12/2023"
$ws.Cells.Item(60, 4).Value = 8.6

$ws.Cells.Item(61, 1).Value = "Nordeste"
$ws.Cells.Item(61, 2).Value = "Outra condição"
$ws.Cells.Item(61, 3).Value = "31/12/2023"
$ws.Cells.Item(61, 4).Value = 0.2

$ws.Cells.Item(62, 1).Value = "Sergipe"
$ws.Cells.Item(62, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(62, 3).Value = "31/12/2016"
$ws.Cells.Item(62, 4).Value = 72.8

$ws.Cells.Item(63, 1).Value = "Sergipe"
$ws.Cells.Item(63, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(63, 3).Value = "31/12/2016"
$ws.Cells.Item(63, 4).Value = 3.5

$ws.Cells.Item(64, 1).Value = "Sergipe"
$ws.Cells.Item(64, 2).Value = "Alugado"
$ws.Cells.Item(64, 3).Value = "31/12/2016"
$ws.Cells.Item(64, 4).Value = 17.5

$ws.Cells.Item(65, 1).Value = "Sergipe"
$ws.Cells.Item(65, 2).Value = "Cedido"
$ws.Cells.Item(65, 3).Value = "31/12/2016"
$ws.Cells.Item(65, 4).Value = 6

$ws.Cells.Item(66, 1).Value = "Sergipe"
$ws.Cells.Item(66, 2).Value = "Outra condição"
$ws.Cells.Item(66, 3).Value = "31/12/2016"
$ws.Cells.Item(66, 4).Value = 0.1

$ws.Cells.Item(67, 1).Value = "Sergipe"
$ws.Cells.Item(67, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(67, 3).Value = "31/12/2017"
$ws.Cells.Item(67, 4).Value = 71.3

$ws.Cells.Item(68, 1).Value = "Sergipe"
$ws.Cells.Item(68, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(68, 3).Value = "31/12/2017"
$ws.Cells.Item(68, 4).Value = 2.8

$ws.Cells.Item(69, 1).Value = "Sergipe"
$ws.Cells.Item(69, 2).Value = "Alugado"
$ws.Cells.Item(69, 3).Value = "31/12/2017"
$ws.Cells.Item(69, 4).Value = 18.6

$ws.Cells.Item(70, 1).Value = "Sergipe"
$ws.Cells.Item(70, 2).Value = "Cedido"
$ws.Cells.Item(70, 3).Value = "31/12/2017"
$ws.Cells.Item(70, 4).Value = 7.1

$ws.Cells.Item(71, 1).Value = "Sergipe"
$ws.Cells.Item(71, 2).Value = "Outra condição"
$ws.Cells.Item(71, 3).Value = "31/12/2017"
$ws.Cells.Item(71, 4).Value = 0.2

$ws.Cells.Item(72, 1).Value = "Sergipe"
$ws.Cells.Item(72, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(72, 3).Value = "31/12/2018"
$ws.Cells.Item(72, 4).Value = 68.59999999999999

$ws.Cells.Item(73, 1).Value = "Sergipe"
$ws.Cells.Item(73, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(73, 3).Value = "31/12/2018"
$ws.Cells.Item(73, 4).Value = 4.9

$ws.Cells.Item(74, 1).Value = "Sergipe"
$ws.Cells.Item(74, 2).Value = "Alugado"
$ws.Cells.Item(74, 3).Value = "31/12/2018"
$ws.Cells.Item(74, 4).Value = 17

$ws.Cells.Item(75, 1).Value = "Sergipe"
$ws.Cells.Item(75, 2).Value = "Cedido"
$ws.Cells.Item(75, 3).Value = "31/12/2018"
$ws.Cells.Item(75, 4).Value = 8.9

$ws.Cells.Item(76, 1).Value = "Sergipe"
$ws.Cells.Item(76, 2).Value = "Outra condição"
$ws.Cells.Item(76, 3).Value = "31/12/2018"
$ws.Cells.Item(76, 4).Value = 0.6

$ws.Cells.Item(77, 1).Value = "Sergipe"
$ws.Cells.Item(77, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(77, 3).Value = "31/12/2019"
$ws.Cells.Item(77, 4).Value = 64.09999999999999

$ws.Cells.Item(78, 1).Value = "Sergipe"
$ws.Cells.Item(78, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(78, 3).Value = "31/12/2019"
$ws.Cells.Item(78, 4).Value = 3.9

$ws.Cells.Item(79, 1).Value = "Sergipe"
$ws.Cells.Item(79, 2).Value = "Alugado"
$ws.Cells.Item(79, 3).Value = "31/12/2019"
$ws.Cells.Item(79, 4).Value = 20.1

$ws.Cells.Item(80, 1).Value = "Sergipe"
$ws.Cells.Item(80, 2).Value = "Cedido"
$ws.Cells.Item(80, 3).Value = "31/12/2019"
$ws.Cells.Item(80, 4).Value = 11.2

$ws.Cells.Item(81, 1).Value = "Sergipe"
$ws.Cells.Item(81, 2).Value = "Outra condição"
$ws.Cells.Item(81, 3).Value = "31/12/2019"
$ws.Cells.Item(81, 4).Value = 0.6

$ws.Cells.Item(82, 1).Value = "Sergipe"
$ws.Cells.Item(82, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(82, 3).Value = "31/12/2022"
$ws.Cells.Item(82, 4).Value = 61.8

$ws.Cells.Item(83, 1).Value = "Sergipe"
$ws.Cells.Item(83, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(83, 3).Value = "31/12/2022"
$ws.Cells.Item(83, 4).Value = 4.8

$ws.Cells.Item(84, 1).Value = "Sergipe"
$ws.Cells.Item(84, 2).Value = "Alugado"
$ws.Cells.Item(84, 3).Value = "31/12/2022"
$ws.Cells.Item(84, 4).Value = 22

$ws.Cells.Item(85, 1).Value = "Sergipe"
$ws.Cells.Item(85, 2).Value = "Cedido"
$ws.Cells.Item(85, 3).Value = "31/12/2022"
$ws.Cells.Item(85, 4).Value = 10.9

$ws.Cells.Item(86, 1).Value = "Sergipe"
$ws.Cells.Item(86, 2).Value = "Outra condição"
$ws.Cells.Item(86, 3).Value = "31/12/2022"
$ws.Cells.Item(86, 4).Value = 0.5

$ws.Cells.Item(87, 1).Value = "Sergipe"
$ws.Cells.Item(87, 2).Value = "Próprio de algum morador - já pago"
$ws.Cells.Item(87, 3).Value = "31/12/2023"
$ws.Cells.Item(87, 4).Value = 64.40000000000001

$ws.Cells.Item(88, 1).Value = "Sergipe"
$ws.Cells.Item(88, 2).Value = "Próprio de algum morador - ainda pagando"
$ws.Cells.Item(88, 3).Value = "31/12/2023"
$ws.Cells.Item(88, 4).Value = 4.1

$ws.Cells.Item(89, 1).Value = "Sergipe"
$ws.Cells.Item(89, 2).Value = "Alugado"
$ws.Cells.Item(89, 3).Value = "31/12/2023"
$ws.Cells.Item(89, 4).Value = 21.1

$ws.Cells.Item(90, 1).Value = "Sergipe"
$ws.Cells.Item(90, 2).Value = "Cedido"
$ws.Cells.Item(90, 3).Value = "31/12/2023"
$ws.Cells.Item(90, 4).Value = 9.699999999999999

$ws.Cells.Item(91, 1).Value = "Sergipe"
$ws.Cells.Item(91, 2).Value = "Outra condição"
$ws.Cells.Item(91, 3).Value = "31/12/2023"
$ws.Cells.Item(91, 4).Value = 0.6
